$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# A new "2020" column (Q) is appended to the table, mirroring the
# formatting of the previous "2019" column (P) for both the header
# row (row 4) and the data row (row 5).
$ws.Range("Q4").Value = 2020
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial($xlPasteFormats)

$ws.Range("Q5").Value = 3.3
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# Reflect the new selection left behind after entering the data.
$ws.Range("R4").Select()
